$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "C4"   = -12.2602
    "B8"   = 5.3137
    "B10"  = 6.635
    "C11"  = -14.02860000000001
    "B12"  = 6.177
    "C12"  = -12.2286
    "C15"  = -13.72289999999999
    "C17"  = -14.35069999999999
    "B18"  = 6.685099999999997
    "B25"  = 6.076299999999996
    "C26"  = -12.55720000000001
    "C27"  = -13.4293
    "C28"  = -13.5971
    "C32"  = -13.48350000000001
    "B37"  = 8.615399999999999
    "C37"  = -12.9188
    "C41"  = -12.48790000000001
    "C47"  = -12.2722
    "C51"  = -11.7337
    "B55"  = 5.456399999999999
    "C65"  = -12.5964
    "B68"  = 4.668699999999994
    "C73"  = -11.21290000000001
    "B77"  = 9.277100000000001
    "B78"  = 9.187599999999993
    "B79"  = 9.116500000000004
    "B80"  = 9.221599999999997
    "B81"  = 6.435500000000006
    "B82"  = 5.813199999999999
    "B84"  = 6.629300000000006
    "C84"  = -12.2945
    "C85"  = -13.17310000000001
    "C89"  = -14.31109999999999
    "C93"  = -10.23499999999999
    "C95"  = -11.97040000000001
    "C98"  = -13.06720000000001
    "C99"  = -12.3523
    "B101" = 5.3203
    "C101" = -13.7341
    "B102" = 7.158399999999999
    "C102" = -12.19350000000001
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
